$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellUpdates = @(
    @{ Cell = 'D2'; Value = '28.908.35' }
    @{ Cell = 'E2'; Value = '  -2.02%  ' }
    @{ Cell = 'D3'; Value = '1.898.45' }
    @{ Cell = 'E3'; Value = '  -3.90%  ' }
    @{ Cell = 'E4'; Value = '  -0.06%  ' }
    @{ Cell = 'D5'; Value = '''324.33' }
    @{ Cell = 'E5'; Value = '  -0.85%  ' }
    @{ Cell = 'E6'; Value = '  -0.16%  ' }
    @{ Cell = 'E7'; Value = '  -1.66%  ' }
    @{ Cell = 'D8'; Value = '''0.3814' }
    @{ Cell = 'E8'; Value = '  -2.53%  ' }
    @{ Cell = 'D9'; Value = '''0.07713' }
    @{ Cell = 'E9'; Value = '  -2.90%  ' }
    @{ Cell = 'D10'; Value = '''0.9751' }
    @{ Cell = 'E10'; Value = '  -1.74%  ' }
    @{ Cell = 'D11'; Value = '''22.03' }
    @{ Cell = 'E11'; Value = '  -3.52%  ' }
    @{ Cell = 'D12'; Value = '1.900.47' }
    @{ Cell = 'E12'; Value = '  -3.99%  ' }
    @{ Cell = 'D13'; Value = '''6.927' }
    @{ Cell = 'E13'; Value = '  -3.59%  ' }
    @{ Cell = 'D14'; Value = '''5.636' }
    @{ Cell = 'E14'; Value = '  -3.72%  ' }
    @{ Cell = 'D15'; Value = '''0.07022' }
    @{ Cell = 'E15'; Value = '  -0.90%  ' }
    @{ Cell = 'E16'; Value = '  -0.13%  ' }
    @{ Cell = 'E17'; Value = '  -4.78%  ' }
    @{ Cell = 'D18'; Value = '''0.000009451' }
    @{ Cell = 'E18'; Value = '  -4.83%  ' }
    @{ Cell = 'D19'; Value = '''16.61' }
    @{ Cell = 'E19'; Value = '  -3.89%  ' }
    @{ Cell = 'E20'; Value = '  -0.13%  ' }
    @{ Cell = 'D21'; Value = '28.865.10' }
    @{ Cell = 'E21'; Value = '  -2.19%  ' }
    @{ Cell = 'D22'; Value = '''5.287' }
    @{ Cell = 'E22'; Value = '  -4.79%  ' }
    @{ Cell = 'D23'; Value = '''10.84' }
    @{ Cell = 'E23'; Value = '  -3.00%  ' }
    @{ Cell = 'D24'; Value = '''2.092' }
    @{ Cell = 'E24'; Value = '  -0.76%  ' }
    @{ Cell = 'D25'; Value = '''157.87' }
    @{ Cell = 'E25'; Value = '  -0.39%  ' }
    @{ Cell = 'D26'; Value = '''18.96' }
    @{ Cell = 'E26'; Value = '  -3.06%  ' }
    @{ Cell = 'D27'; Value = '''5.618' }
    @{ Cell = 'E27'; Value = '  -2.96%  ' }
    @{ Cell = 'D28'; Value = '''117.21' }
    @{ Cell = 'E28'; Value = '  -1.96%  ' }
    @{ Cell = 'D29'; Value = '''1.830' }
    @{ Cell = 'E29'; Value = '  -4.11%  ' }
    @{ Cell = 'D30'; Value = '''0.09238' }
    @{ Cell = 'E30'; Value = '  -1.97%  ' }
    @{ Cell = 'D31'; Value = '''0.8598' }
    @{ Cell = 'E31'; Value = '  -3.78%  ' }
    @{ Cell = 'D32'; Value = '''5.078' }
    @{ Cell = 'E32'; Value = '  -2.99%  ' }
    @{ Cell = 'D33'; Value = '''1.238' }
    @{ Cell = 'E33'; Value = '  -6.35%  ' }
    @{ Cell = 'D34'; Value = '''2.997' }
    @{ Cell = 'E34'; Value = '  -6.01%  ' }
    @{ Cell = 'D35'; Value = '''0.05683' }
    @{ Cell = 'E35'; Value = '  -2.32%  ' }
    @{ Cell = 'D36'; Value = '''1.140' }
    @{ Cell = 'E36'; Value = '  -2.56%  ' }
    @{ Cell = 'E37'; Value = '  -0.03%  ' }
    @{ Cell = 'D38'; Value = '''0.02026' }
    @{ Cell = 'E38'; Value = '  -3.58%  ' }
    @{ Cell = 'D39'; Value = '''0.5473' }
    @{ Cell = 'E39'; Value = '  -4.34%  ' }
    @{ Cell = 'D40'; Value = '''7.375' }
    @{ Cell = 'E40'; Value = '  -5.22%  ' }
    @{ Cell = 'E41'; Value = '  -3.02%  ' }
    @{ Cell = 'D42'; Value = '''9.248' }
    @{ Cell = 'E42'; Value = '  -4.19%  ' }
    @{ Cell = 'D43'; Value = '''2.759' }
    @{ Cell = 'E43'; Value = '  +0.03%  ' }
    @{ Cell = 'D44'; Value = '''0.5148' }
    @{ Cell = 'E44'; Value = '  -4.02%  ' }
    @{ Cell = 'D45'; Value = '''11.18' }
    @{ Cell = 'E45'; Value = '  -4.78%  ' }
    @{ Cell = 'D46'; Value = '''0.06805' }
    @{ Cell = 'E46'; Value = '  -1.72%  ' }
    @{ Cell = 'B47'; Value = 'RenderToken' }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' }
    @{ Cell = 'D47'; Value = '''2.067' }
    @{ Cell = 'E47'; Value = '  -5.52%  ' }
    @{ Cell = 'B48'; Value = 'PEPE' }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe' }
    @{ Cell = 'D48'; Value = '''0.000002613' }
    @{ Cell = 'E48'; Value = '  -17.42%  ' }
    @{ Cell = 'D49'; Value = '''110.17' }
    @{ Cell = 'E49'; Value = '  -3.45%  ' }
    @{ Cell = 'D50'; Value = '''1.768' }
    @{ Cell = 'E50'; Value = '  -3.36%  ' }
    @{ Cell = 'E51'; Value = '  -0.16%  ' }
)

foreach ($update in $cellUpdates) {
    $ws.Range($update.Cell).Value = $update.Value
}
